$d = $word.ActiveDocument

function Replace-WithXml {
    param($SearchText, $InnerXml)
    $r = $d.Content
    $found = $r.Find.Execute($SearchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $SearchText"
    }
    # Re-materialize the found span as a fresh Range object; InsertXML on the
    # Find's own range object behaves as an insert (duplicating text) after
    # earlier document mutations, whereas a freshly constructed Range with
    # the same Start/End correctly replaces its content.
    $target = $d.Range($r.Start, $r.End)
    $frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/dummy.xml" pkg:contentType="x"><pkg:xmlData><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $InnerXml + '</w:p></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($frag)
}

# 1. Heading: "Bug Test – Fine Calc Incorrect" -> two runs
Replace-WithXml "Bug Test – Fine Calc Incorrect" '<w:r><w:t xml:space="preserve">Bug Test – </w:t></w:r><w:r><w:t>No Fine for one day overdue</w:t></w:r>'

# 2. " ...use case, etc…" -> split off "etc" with proofErr
Replace-WithXml " (aka test set or test suite) are a set of test scripts that cover a specific functional area, business process, use case, etc…" '<w:r><w:t xml:space="preserve"> (aka test set or test suite) are a set of test scripts that cover a specific functional area, business process, use case, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>etc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>…</w:t></w:r>'

# 3. "For example, ... exception, etc…" -> split off "etc" with proofErr
Replace-WithXml "For example, a given test scenario might cover an order placed on a web site by an existing customer (another scenario might cover orders placed by new customers) – test scripts within the scenario might cover a single item order, a multiple items order, quantity not on hand exception, etc…" '<w:r><w:t xml:space="preserve">For example, a given test scenario might cover an order placed on a web site by an existing customer (another scenario might cover orders placed by new customers) – test scripts within the scenario might cover a single item order, a multiple items order, quantity not on hand exception, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>etc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>…</w:t></w:r>'

# 4. "1.2 returnBookAutomatedTest" -> split with proofErr around returnBookAutomatedTest
Replace-WithXml "1.2 returnBookAutomatedTest" '<w:r><w:t xml:space="preserve">1.2 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>returnBookAutomatedTest</w:t></w:r><w:proofErr w:type="spellEnd"/>'

# 5. "Create a book, title t, author a, cNo c1" -> split with proofErr around cNo
Replace-WithXml "Create a book, title t, author a, cNo c1" '<w:r><w:t xml:space="preserve">Create a book, title t, author a, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>cNo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> c1</w:t></w:r>'

# 6. "TimothyPickard" in table cell -> wrap with proofErr
Replace-WithXml "TimothyPickard" '<w:proofErr w:type="spellStart"/><w:r><w:t>TimothyPickard</w:t></w:r><w:proofErr w:type="spellEnd"/>'
